$d = $word.ActiveDocument

# --- Edit 1: "This document will guide you through ..." paragraph ---
# Replace the old two-sentence tail with the new wording (kept as one
# contiguous run of matching run-properties; downstream runs describing
# "One of the / examples / is generated ... / ." share identical
# sz=24/szCs=24 formatting with the text that precedes them).
$d.Content.Find.Execute(
    "There are two examples. Their levels are the same, but one generates C++ and the other XML output.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "One of the examples is generated through a code described in this guide.",
    2
) | Out-Null

# --- Edit 2: "You can compare it to the output from the example project. ..." paragraph ---
$d.Content.Find.Execute(
    "There is another example with the same level, but which generates an XML file.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There are other examples with the same level, but which generate different formats.",
    2
) | Out-Null
